$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1858.2407
$ws.Range("I138").Value = 1686.3541
$ws.Range("J138").Value = 3233.3333
$ws.Range("K138").Value = 5059.0623
$ws.Range("L138").Value = 9699.999899999999
$ws.Range("M138").Value = 80.9377000000004
$ws.Range("N138").Value = -19979.9999

$ws.Range("H141").Value = 945.125
$ws.Range("I141").Value = 625.76
$ws.Range("J141").Value = 2085.7144
$ws.Range("K141").Value = 1877.28
$ws.Range("L141").Value = 6257.1432
$ws.Range("M141").Value = 3302.72
$ws.Range("N141").Value = -16617.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8209.864
$ws.Range("I32").Value = 4316.8623
$ws.Range("J32").Value = 45842.223
$ws.Range("K32").Value = 4316.8623
$ws.Range("L32").Value = 45842.223
$ws.Range("M32").Value = -4029.8623
$ws.Range("N32").Value = -46416.223

$ws.Range("H37").Value = 8950
$ws.Range("I37").Value = 3750
$ws.Range("J37").Value = 10683.333
$ws.Range("K37").Value = 3750
$ws.Range("L37").Value = 10683.333
$ws.Range("M37").Value = -3477
$ws.Range("N37").Value = -11229.333

$ws.Range("H55").Value = 17564.715
$ws.Range("J55").Value = 17564.715
$ws.Range("L55").Value = 17564.715
$ws.Range("N55").Value = -18194.715

$ws.Range("H132").Value = 18420.543
$ws.Range("I132").Value = 1410.625
$ws.Range("K132").Value = 4231.875
$ws.Range("M132").Value = -1701.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H34").Value = 10000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 10000
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 10000
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -10228

$ws.Range("H82").Value = 15317.941
$ws.Range("I82").Value = 8333.75
$ws.Range("J82").Value = 21526.111
$ws.Range("K82").Value = 8333.75
$ws.Range("L82").Value = 21526.111
$ws.Range("M82").Value = -7950.75
$ws.Range("N82").Value = -22292.111

$ws.Range("H85").Value = 15317.941
$ws.Range("I85").Value = 8333.75
$ws.Range("J85").Value = 21526.111
$ws.Range("K85").Value = 8333.75
$ws.Range("L85").Value = 21526.111
$ws.Range("M85").Value = -7007.75
$ws.Range("N85").Value = -24178.111

$ws.Range("H134").Value = 17910.299
$ws.Range("I134").Value = 19912.033
$ws.Range("J134").Value = 752.5714
$ws.Range("K134").Value = 59736.099
$ws.Range("L134").Value = 2257.7142
$ws.Range("M134").Value = -57201.099
$ws.Range("N134").Value = -7327.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 10447.75
$ws.Range("I50").Value = 7500
$ws.Range("J50").Value = 11430.333
$ws.Range("K50").Value = 7500
$ws.Range("L50").Value = 11430.333
$ws.Range("M50").Value = -6875
$ws.Range("N50").Value = -12680.333

$ws.Range("H51").Value = 11636.363
$ws.Range("J51").Value = 12455.556
$ws.Range("L51").Value = 12455.556
$ws.Range("N51").Value = -13927.556

$ws.Range("H60").Value = 8890.5
$ws.Range("I60").Value = 5599.6665
$ws.Range("J60").Value = 10300.857
$ws.Range("K60").Value = 5599.6665
$ws.Range("L60").Value = 10300.857
$ws.Range("M60").Value = -5088.6665
$ws.Range("N60").Value = -11322.857

$ws.Range("H61").Value = 11636.363
$ws.Range("J61").Value = 12455.556
$ws.Range("L61").Value = 12455.556
$ws.Range("N61").Value = -13151.556

$ws.Range("H68").Value = 15041
$ws.Range("I68").Value = 8892
$ws.Range("J68").Value = 19960.2
$ws.Range("K68").Value = 8892
$ws.Range("L68").Value = 19960.2
$ws.Range("M68").Value = -8143
$ws.Range("N68").Value = -21458.2

$ws.Range("H71").Value = 15041
$ws.Range("I71").Value = 8892
$ws.Range("J71").Value = 19960.2
$ws.Range("K71").Value = 26676
$ws.Range("L71").Value = 59880.60000000001
$ws.Range("M71").Value = -22932
$ws.Range("N71").Value = -67368.60000000001

$ws.Range("H74").Value = 15314.5
$ws.Range("I74").Value = 5800
$ws.Range("J74").Value = 16673.715
$ws.Range("K74").Value = 5800
$ws.Range("L74").Value = 16673.715
$ws.Range("M74").Value = -4926
$ws.Range("N74").Value = -18421.715

$ws.Range("H77").Value = 15314.5
$ws.Range("I77").Value = 5800
$ws.Range("J77").Value = 16673.715
$ws.Range("K77").Value = 17400
$ws.Range("L77").Value = 50021.145
$ws.Range("M77").Value = -13032
$ws.Range("N77").Value = -58757.145

$ws.Range("H88").Value = 28894
$ws.Range("J88").Value = 28894
$ws.Range("L88").Value = 28894
$ws.Range("N88").Value = -29706

$ws.Range("H91").Value = 28894
$ws.Range("J91").Value = 28894
$ws.Range("L91").Value = 28894
$ws.Range("N91").Value = -31702

$ws.Range("H109").Value = 42250
$ws.Range("J109").Value = 42250
$ws.Range("L109").Value = 42250
$ws.Range("N109").Value = -44330

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 794.7458
$ws.Range("I5").Value = 487.72726
$ws.Range("J5").Value = 865.1042
$ws.Range("K5").Value = 1463.18178
$ws.Range("L5").Value = 2595.3126
$ws.Range("M5").Value = -1351.18178
$ws.Range("N5").Value = -2819.3126

$ws.Range("H39").Value = 10022.667
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 10022.667
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 30068.001
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -30656.001

$ws.Range("H40").Value = 111.5
$ws.Range("I40").Value = 111.5
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 446
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -377
$ws.Range("N40").ClearContents()

$ws.Range("H135").Value = 794.7458
$ws.Range("I135").Value = 487.72726
$ws.Range("J135").Value = 865.1042
$ws.Range("K135").Value = 4389.54534
$ws.Range("L135").Value = 7785.9378
$ws.Range("M135").Value = -1854.54534
$ws.Range("N135").Value = -12855.9378

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3875.35
$ws.Range("I70").Value = 3634.75
$ws.Range("J70").Value = 4236.25
$ws.Range("K70").Value = 3634.75
$ws.Range("L70").Value = 4236.25
$ws.Range("M70").Value = -3364.75
$ws.Range("N70").Value = -4776.25

$ws.Range("H73").Value = 3875.35
$ws.Range("I73").Value = 3634.75
$ws.Range("J73").Value = 4236.25
$ws.Range("K73").Value = 3634.75
$ws.Range("L73").Value = 4236.25
$ws.Range("M73").Value = -2698.75
$ws.Range("N73").Value = -6108.25

$ws.Range("H132").Value = 4675.3335
$ws.Range("I132").Value = 4610.4
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 13831.2
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -11301.2
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1090.5714
$ws.Range("I136").Value = 755.15625
$ws.Range("J136").Value = 4668.3335
$ws.Range("K136").Value = 2265.46875
$ws.Range("L136").Value = 14005.0005
$ws.Range("M136").Value = 284.53125
$ws.Range("N136").Value = -19105.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
